$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 511, shifting rows 511:609 down to 512:610.
$ws.Rows.Item(511).Insert()

# Populate the newly inserted row 511 with the new weekly price entry.
$ws.Cells.Item(511, 1).Value = 5
$ws.Cells.Item(511, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(511, 3).Value = "Maule"
$ws.Cells.Item(511, 4).Value = 45209
$ws.Cells.Item(511, 5).Value = 7
$ws.Cells.Item(511, 6).Value = 100114013
$ws.Cells.Item(511, 7).Value = "Zanahoria"
$ws.Cells.Item(511, 8).Value = "Sin especificar"
$ws.Cells.Item(511, 9).Value = "Primera"
$ws.Cells.Item(511, 10).Value = 500
$ws.Cells.Item(511, 11).Value = 5000
$ws.Cells.Item(511, 12).Value = 5000
$ws.Cells.Item(511, 13).Value = 5000
$ws.Cells.Item(511, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(511, 15).Value = "Región de Ñuble"
$ws.Cells.Item(511, 16).Value = 250
$ws.Cells.Item(511, 17).Value = 20
$ws.Cells.Item(511, 18).Value = "Hortaliza"
